$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(1).Insert()
$ws.Range("A17").Value = "DET -> this | that | these | those"
$ws.Range("A5").Select()
